$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I1, J1) re-using the same header style ("s=1")
# that B1..H1 already carry, by copy/pasting the formatting from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 / IF values for every data row (2-41).
$data = @{
    2  = @(7, 7)
    3  = @(1, 1)
    4  = @(8, 8)
    5  = @(6, 6)
    6  = @(6, 6)
    7  = @(9, 9)
    8  = @(10, 10)
    9  = @(7, 7)
    10 = @(6, 6)
    11 = @(6, 6)
    12 = @(8, 8)
    13 = @(7, 7)
    14 = @(7, 7)
    15 = @(8, 9)
    16 = @(7, 7)
    17 = @(1, 2)
    18 = @(7, 7)
    19 = @(4, 5)
    20 = @(1, 2)
    21 = @(5, 6)
    22 = @(4, 4)
    23 = @(8, 9)
    24 = @(7, 8)
    25 = @(6, 7)
    26 = @(5, 5)
    27 = @(5, 6)
    28 = @(1, 2)
    29 = @(1, 3)
    30 = @(6, 6)
    31 = @(4, 5)
    32 = @(5, 6)
    33 = @(9, 9)
    34 = @(6, 6)
    35 = @(5, 6)
    36 = @(7, 7)
    37 = @(9, 9)
    38 = @(7, 7)
    39 = @(9, 9)
    40 = @(1, 2)
    41 = @(8, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
